$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # 展览
$ws1.Range("F2").Value = 224
$ws1.Range("F4").Value = 19723
$ws1.Range("F5").Value = 795
$ws1.Range("F6").Value = 0
$ws1.Range("F7").Value = 1092
$ws1.Range("F8").Value = 0
$ws1.Range("F9").Value = 7468
$ws1.Range("F10").Value = 497
$ws1.Range("F12").Value = 257
$ws1.Range("F15").Value = 107
$ws1.Range("F16").Value = 5
$ws1.Range("F17").Value = 232
$ws1.Range("F19").Value = 0
$ws1.Range("F20").Value = 388
$ws1.Range("F23").Value = 0
$ws1.Range("F24").Value = 0
$ws1.Range("F27").Value = 0
$ws1.Range("F30").Value = 172
$ws1.Range("F31").Value = 0
$ws1.Range("F32").Value = 558
$ws1.Range("F33").Value = 0
$ws1.Range("F34").Value = 0
$ws1.Range("F35").Value = 24
$ws1.Range("F37").Value = 22
$ws1.Range("F38").Value = 12551
$ws1.Range("F39").Value = 0
$ws1.Range("F40").Value = 65
$ws1.Range("F42").Value = 0
$ws1.Range("F43").Value = 254
$ws1.Range("F44").Value = 350
$ws1.Range("F45").Value = 3981
$ws1.Range("F47").Value = 0

$ws4 = $wb.Worksheets.Item(4)   # 全部类型
$ws4.Range("F2").Value = 224
$ws4.Range("F4").Value = 0
$ws4.Range("F5").Value = 795
$ws4.Range("F6").Value = 306
$ws4.Range("F9").Value = 7469
$ws4.Range("F10").Value = 0
$ws4.Range("F12").Value = 0
$ws4.Range("F14").Value = 149
$ws4.Range("F15").Value = 107
$ws4.Range("F17").Value = 0
$ws4.Range("F18").Value = 188
$ws4.Range("F20").Value = 0
$ws4.Range("F21").Value = 70
$ws4.Range("F23").Value = 0
$ws4.Range("F24").Value = 51
$ws4.Range("F26").Value = 315
$ws4.Range("F28").Value = 25
$ws4.Range("F29").Value = 0
$ws4.Range("F31").Value = 5219
$ws4.Range("F32").Value = 558
$ws4.Range("F34").Value = 52
$ws4.Range("F35").Value = 0
$ws4.Range("F36").Value = 2792
$ws4.Range("F38").Value = 86
$ws4.Range("F39").Value = 22
$ws4.Range("F40").Value = 12551
$ws4.Range("F41").Value = 0
$ws4.Range("F42").Value = 65
$ws4.Range("F44").Value = 54
$ws4.Range("F46").Value = 0
$ws4.Range("F47").Value = 3981
$ws4.Range("F49").Value = 0
